$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.00666173205033309
$ws.Range("C2").Value = 0.0111028867505551
$ws.Range("D2").Value = 0.00296076980014804
$ws.Range("E2").Value = 0.00148038490007402
$ws.Range("F2").Value = 0.0111028867505551
$ws.Range("G2").Value = 0.117690599555885
$ws.Range("H2").Value = 0.0399703923019985
$ws.Range("I2").Value = 0.0162842339008142
$ws.Range("J2").Value = 0.00370096225018505
$ws.Range("K2").Value = 0.0125832716506292
$ws.Range("L2").Value = 0.00518134715025907
$ws.Range("M2").Value = 0.00370096225018505
$ws.Range("N2").Value = 0.00370096225018505
$ws.Range("O2").Value = 0.0185048112509252
$ws.Range("P2").Value = 0.00148038490007402
$ws.Range("Q2").Value = 0.00666173205033309
$ws.Range("R2").Value = 0.977054034048853
$ws.Range("S2").Value = 0.0074019245003701
$ws.Range("T2").Value = 0.0429311621021466
$ws.Range("U2").Value = 0.0266469282013323
$ws.Range("V2").Value = 0.00814211695040711
$ws.Range("W2").Value = 0.00148038490007402
$ws.Range("X2").Value = 0.0125832716506292
$ws.Range("B3").Value = 0.00222057735011103
$ws.Range("C3").Value = 0.686158401184308
$ws.Range("D3").Value = 0.00222057735011103
$ws.Range("E3").Value = 0.0421909696521095
$ws.Range("F3").Value = 0.00074019245003701
$ws.Range("G3").Value = 0.00888230940044412
$ws.Range("H3").Value = 0.920059215396003
$ws.Range("I3").Value = 0.973353071798668
$ws.Range("J3").Value = 0.00074019245003701
$ws.Range("K3").Value = 0.00074019245003701
$ws.Range("L3").Value = 0.00666173205033309
$ws.Range("M3").Value = 0.00592153960029608
$ws.Range("N3").Value = 0.974833456698742
$ws.Range("O3").Value = 0.00074019245003701
$ws.Range("R3").Value = 0.00148038490007402
$ws.Range("S3").Value = 0.982975573649149
$ws.Range("T3").Value = 0.0895632864544782
$ws.Range("U3").Value = 0.00518134715025907
$ws.Range("W3").Value = 0.00962250185048113
$ws.Range("X3").Value = 0.00370096225018505
$ws.Range("B4").Value = 0.990377498149519
$ws.Range("C4").Value = 0.00370096225018505
$ws.Range("D4").Value = 0.988156920799408
$ws.Range("E4").Value = 0.0133234641006662
$ws.Range("F4").Value = 0.982235381199112
$ws.Range("G4").Value = 0.86380458919319
$ws.Range("H4").Value = 0.00296076980014804
$ws.Range("I4").Value = 0.00074019245003701
$ws.Range("J4").Value = 0.0259067357512953
$ws.Range("K4").Value = 0.981495188749075
$ws.Range("L4").Value = 0.987416728349371
$ws.Range("M4").Value = 0.988156920799408
$ws.Range("N4").Value = 0.00222057735011103
$ws.Range("O4").Value = 0.978534418948927
$ws.Range("P4").Value = 0.998519615099926
$ws.Range("Q4").Value = 0.99259807549963
$ws.Range("R4").Value = 0.0162842339008142
$ws.Range("S4").Value = 0.00074019245003701
$ws.Range("T4").Value = 0.0259067357512953
$ws.Range("U4").Value = 0.962250185048113
$ws.Range("V4").Value = 0.983715766099186
$ws.Range("W4").Value = 0.987416728349371
$ws.Range("X4").Value = 0.97779422649889
$ws.Range("B5").Value = 0.00074019245003701
$ws.Range("C5").Value = 0.299037749814952
$ws.Range("D5").Value = 0.00666173205033309
$ws.Range("E5").Value = 0.94300518134715
$ws.Range("F5").Value = 0.00592153960029608
$ws.Range("G5").Value = 0.00888230940044412
$ws.Range("H5").Value = 0.0362694300518135
$ws.Range("I5").Value = 0.00962250185048113
$ws.Range("J5").Value = 0.969652109548483
$ws.Range("K5").Value = 0.00444115470022206
$ws.Range("L5").Value = 0.00074019245003701
$ws.Range("M5").Value = 0.00074019245003701
$ws.Range("N5").Value = 0.0192450037009623
$ws.Range("O5").Value = 0.00148038490007402
$ws.Range("Q5").Value = 0.00074019245003701
$ws.Range("R5").Value = 0.00518134715025907
$ws.Range("S5").Value = 0.00888230940044412
$ws.Range("T5").Value = 0.84159881569208
$ws.Range("U5").Value = 0.00592153960029608
$ws.Range("X5").Value = 0.00444115470022206
